$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ r=27; A="2022030211"; B="22"; D="62"; E="东风1级" },
    @{ r=28; A="2022030210"; B="21"; C="59"; D="66"; E="东北风1级" },
    @{ r=29; A="2022030209"; B="19"; C="58"; D="72"; E="东风1级" },
    @{ r=30; A="2022030208"; B="18"; C="58"; D="77"; E="东风1级" },
    @{ r=31; A="2022030207"; B="17"; C="60"; D="84"; E="东北风1级" },
    @{ r=32; A="2022030206"; B="17"; C="62"; D="85"; E="东南风1级" },
    @{ r=33; A="2022030205"; B="18"; C="63"; D="85"; E="东风1级" },
    @{ r=34; A="2022030204"; B="18"; C="65"; D="91"; E="南风1级" },
    @{ r=35; A="2022030203"; B="18"; C="67"; D="96"; E="南风1级" },
    @{ r=36; A="2022030202"; B="18"; C="71"; D="97"; E="东风0级" },
    @{ r=37; A="2022030201"; B="18"; C="74"; D="98"; E="东风0级" },
    @{ r=38; A="2022030200"; B="18"; C="73"; D="98"; E="东北风1级" },
    @{ r=39; A="2022030223"; B="19"; C="69"; D="95"; E="南风0级" },
    @{ r=40; A="2022030222"; B="19"; C="68"; D="93"; E="东北风1级" },
    @{ r=41; A="2022030221"; B="19"; C="64"; D="89"; E="东风0级" },
    @{ r=42; A="2022030220"; B="20"; C="60"; D="83"; E="北风1级" },
    @{ r=43; A="2022030219"; B="22"; C="51"; D="73"; E="东北风1级" },
    @{ r=44; A="2022030218"; B="24"; C="47"; D="63"; E="西南风1级" },
    @{ r=45; A="2022030217"; B="25"; C="47"; D="56"; E="西南风1级" },
    @{ r=46; A="2022030216"; B="26"; C="47"; D="54"; E="南风1级" },
    @{ r=47; A="2022030215"; B="26"; C="44"; D="53"; E="南风1级" },
    @{ r=48; A="2022030214"; B="25"; C="48"; D="57"; E="西风1级" },
    @{ r=49; A="2022030213"; B="25"; C="53"; D="59"; E="西南风1级" },
    @{ r=50; A="2022030212"; B="24"; C="55"; D="63"; E="西南风1级" },
    @{ r=51; A="2022030211"; B="24"; C="62"; D="68"; E="东南风1级" }
)

foreach ($row in $rows) {
    foreach ($col in @("A","B","C","D","E")) {
        if ($row.ContainsKey($col)) {
            $addr = "$col$($row.r)"
            $ws.Range($addr).NumberFormat = "@"
            $ws.Range($addr).Value = $row[$col]
        }
    }
}